# Add season-record columns (Wins / Losses / Ties) to the player table.
# The team went 91-71-0 in 1996, so every player row gets the same record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting from the last existing header cell (AB1, style
# index 1: bold + bordered + centered/top-aligned) onto the three new header
# cells, then fill in their labels.
$ws.Range("AB1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)

$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Season record is the same for every player row (team-wide record).
$wins = 91
$losses = 71
$ties = 0

for ($row = 2; $row -le 45; $row++) {
    $ws.Cells.Item($row, 29).Value = $wins
    $ws.Cells.Item($row, 30).Value = $losses
    $ws.Cells.Item($row, 31).Value = $ties
}
